# Weekly update: prepend two new price records (dated 45194) for the
# "Fruta, Feria Lagunitas de Puerto Montt - Plátano" sheet.
#
# This is implemented as an insertion of two new rows at 882:883, which
# shifts all the existing data (previously rows 882-966) down to rows
# 884-968, and then populating the two freshly inserted rows with the
# new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 882, pushing the rest of
# the table (rows 882:966) down to rows 884:968.
$ws.Rows("882:883").Insert()

# --- New row 882 ---
$ws.Cells.Item(882, 1).Value = 4
$ws.Cells.Item(882, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(882, 3).Value = "Los Lagos"
$ws.Cells.Item(882, 4).Value = 45194
$ws.Cells.Item(882, 5).Value = 10
$ws.Cells.Item(882, 6).Value = "Fruta"
$ws.Cells.Item(882, 7).Value = 100108
$ws.Cells.Item(882, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(882, 9).Value = 100108006
$ws.Cells.Item(882, 10).Value = "Plátano"
$ws.Cells.Item(882, 11).Value = "Sin especificar"
$ws.Cells.Item(882, 12).Value = "Pintón"
$ws.Cells.Item(882, 13).Value = 300
$ws.Cells.Item(882, 14).Value = 20000
$ws.Cells.Item(882, 15).Value = 20000
$ws.Cells.Item(882, 16).Value = 20000
$ws.Cells.Item(882, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(882, 18).Value = "Ecuador"
$ws.Cells.Item(882, 19).Value = 1000
$ws.Cells.Item(882, 20).Value = 20

# --- New row 883 ---
$ws.Cells.Item(883, 1).Value = 4
$ws.Cells.Item(883, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(883, 3).Value = "Los Lagos"
$ws.Cells.Item(883, 4).Value = 45194
$ws.Cells.Item(883, 5).Value = 10
$ws.Cells.Item(883, 6).Value = "Fruta"
$ws.Cells.Item(883, 7).Value = 100108
$ws.Cells.Item(883, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(883, 9).Value = 100108006
$ws.Cells.Item(883, 10).Value = "Plátano"
$ws.Cells.Item(883, 11).Value = "Sin especificar"
$ws.Cells.Item(883, 12).Value = "Primera Pintón"
$ws.Cells.Item(883, 13).Value = 800
$ws.Cells.Item(883, 14).Value = 22000
$ws.Cells.Item(883, 15).Value = 23000
$ws.Cells.Item(883, 16).Value = 22500
$ws.Cells.Item(883, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(883, 18).Value = "Ecuador"
$ws.Cells.Item(883, 19).Value = 1125
$ws.Cells.Item(883, 20).Value = 20
